# Append a new "." run right after the "Epale" run, with identical
# character formatting (Times New Roman / Bold / color 7030A0 / 24pt),
# keeping it as its own <w:r> rather than merging into the existing
# "Epale" text run.

$d = $word.ActiveDocument

# Locate the "Epale" text and collapse the range to its end point.
$rng = $d.Content
$found = $rng.Find.Execute("Epale", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# Insert the period right after "Epale" (inherits the surrounding
# run's formatting automatically).
$rng.InsertAfter(".")

# Touch the Bold property (off, then back on) so the newly inserted
# text is kept as a distinct run instead of being coalesced into the
# preceding "Epale" run, even though the effective formatting is the
# same in the end.
$rng.Font.Bold = 0
$rng.Font.Bold = 1
